$d = $word.ActiveDocument
$pairs = @(
    @("29-12=", "13+0="),
    @("32-31=", "16-11="),
    @("36+36=", "93-4="),
    @("61-3=", "26+9="),
    @("80-9=", "93-71="),
    @("48+0=", "82+2="),
    @("95-62=", "58+3="),
    @("58-14=", "50-14="),
    @("6+65=", "66+17="),
    @("12+5=", "13+18="),
    @("98-14=", "9+40="),
    @("19+35=", "75+8="),
    @("41+54=", "83+1="),
    @("0+13=", "38-29="),
    @("35-17=", "56-15="),
    @("45-22=", "0+90="),
    @("28+63=", "11+66="),
    @("92-50=", "87-19="),
    @("86-45=", "80-17="),
    @("78-60=", "20+34="),
    @("37-6=", "77+12="),
    @("49+9=", "57-37="),
    @("70-55=", "48+48="),
    @("51+13=", "12+79="),
    @("61+38=", "67-42="),
    @("4+70=", "72-20="),
    @("99-35=", "27-22="),
    @("93-64=", "88-40="),
    @("94-23=", "8+70="),
    @("32-5=", "46-24="),
    @("3+67=", "96-67="),
    @("85-33=", "91-46="),
    @("38+50=", "41+58="),
    @("60-55=", "2+22="),
    @("57-1=", "89+3="),
    @("77-42=", "0+99="),
    @("81+12=", "67-1="),
    @("35-12=", "26+7="),
    @("61-26=", "96+2="),
    @("78-5=", "16-7="),
    @("25-6=", "67-40="),
    @("99-41=", "69-39="),
    @("10+35=", "89-60="),
    @("7+46=", "64-57="),
    @("74+11=", "3+70="),
    @("30+6=", "46-15="),
    @("66-54=", "17+38="),
    @("10+19=", "51-46="),
    @("34+60=", "98-50="),
    @("60+38=", "77-2="),
    @("25+51=", "75-57="),
    @("77-25=", "31+13="),
    @("48-27=", "44-22="),
    @("47-0=", "19+9="),
    @("19+58=", "8+12="),
    @("17-3=", "95-35="),
    @("16+63=", "7+17="),
    @("52-8=", "22+36="),
    @("65-26=", "5+56="),
    @("54-31=", "35+47="),
    @("62-20=", "86-35="),
    @("98-32=", "85-28="),
    @("74-42=", "73-39="),
    @("57-25=", "93-80="),
    @("64+28=", "54-24="),
    @("31-31=", "60+11="),
    @("51+44=", "44+24="),
    @("22-17=", "8+66="),
    @("2+67=", "60+11="),
    @("67-61=", "92-33="),
    @("99-85=", "14-2="),
    @("67-0=", "76-64="),
    @("98-64=", "16+33="),
    @("28-19=", "26+72="),
    @("89-80=", "58+15="),
    @("56-44=", "53+42="),
    @("97-44=", "82-79="),
    @("92+0=", "30-21="),
    @("27-8=", "86-58="),
    @("82-21=", "70-20="),
    @("55-51=", "24+33="),
    @("23+7=", "94+1="),
    @("9+18=", "3+41="),
    @("46-17=", "9+46="),
    @("82-33=", "71+6="),
    @("88-65=", "59-39="),
    @("91-89=", "71+12="),
    @("77-26=", "76-54="),
    @("47-18=", "47-42="),
    @("14+64=", "38-23="),
    @("57+5=", "50+8="),
    @("70-51=", "98-19="),
    @("68+12=", "22+20="),
    @("1+32=", "57+38="),
    @("0+53=", "85-28="),
    @("43+34=", "29+17="),
    @("86-81=", "93-47="),
    @("7+35=", "84-13="),
    @("97-12=", "49-43="),
    @("15+64=", "92-6="),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
